$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell H1 - "Save", using same style as other header cells (copy format from G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data values for the new "Save" column (H2:H12)
$saveValues = @(0, 0, 1, 0, 0, 1, 1, 1, 0, 0, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
